$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.232499999999996
$ws.Range("A3").Value = -21.49020000000003
$ws.Range("B5").Value = 4.9854
$ws.Range("A14").Value = -20.53299999999997
$ws.Range("A21").Value = -21.26320000000001
$ws.Range("A23").Value = -21.41840000000003
$ws.Range("A25").Value = -22.34580000000004
